# "finished 2017 csv file"
# Rename the weekly sheets so the zero-padded numbers 01-09 become 1-9
# (sem10..sem14 already have no leading zero and stay unchanged).
$wb = $excel.ActiveWorkbook

$renames = @{
    "sem01" = "sem1"
    "sem02" = "sem2"
    "sem03" = "sem3"
    "sem04" = "sem4"
    "sem05" = "sem5"
    "sem06" = "sem6"
    "sem07" = "sem7"
    "sem08" = "sem8"
    "sem09" = "sem9"
}

foreach ($oldName in $renames.Keys) {
    $wb.Worksheets.Item($oldName).Name = $renames[$oldName]
}
